# Update the workbook's single data row (row 2) on sheet "Hoja1" to reflect
# the new course ("SIST_7014" - Bases de Datos Geoespaciales) that replaced
# the old one ("EDUC_7094" - Una Nueva Mirada a la Orientación y Asesoramiento
# Familiar), per commit "180625 Código banner incorrecto".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Row 2 cell values -------------------------------------------------
$ws.Range("A2").Value = 77
$ws.Range("E2").Value = "SIST_7014"
$ws.Range("F2").Value = "Bases de Datos Geoespaciales"
$ws.Range("G2").Value = "Pucha Cofrep Franz Leonardo"
$ws.Range("H2").Value = "fapucha@utpl.edu.ec"
$ws.Range("I2").Value = 1104483498
$ws.Range("K2").Value = "Ingenierías y Arquitectura"
$ws.Range("L2").Value = "Especialización en Gestión de Geoinformación con mención en Proyectos de Ingeniería"
$ws.Range("M2").Value = "González Jaramillo Víctor Hugo"
$ws.Range("N2").Value = "vhgonzalez@utpl.edu.ec"
$ws.Range("O2").Value = 1
$ws.Range("Q2").Value = "Total 144: ACD_32 APE_16 AA_96"
$ws.Range("R2").Value = "Unidad de Formación Disciplinar Avanzada"
$ws.Range("S2").Value = 2
$ws.Range("T2").Value = "Especialización"
$ws.Range("V2").Value = "https://utpl.instructure.com/courses/72916"
$ws.Range("X2").Value = "SIST_7014_META"
$ws.Range("Y2").Value = "180625 Código banner incorrecto"
$ws.Range("AA2").Value = "SIST_7014"
$ws.Range("AB2").Value = 45841
$ws.Range("AC2").Value = 45841

# --- Hyperlinks ----------------------------------------------------------
# N2 (director's e-mail) becomes a live mailto: hyperlink; V2 keeps its
# hyperlink but now targets the new Canvas course URL above.
$ws.Hyperlinks.Add($ws.Range("N2"), "mailto:vhgonzalez@utpl.edu.ec")
$ws.Range("N2").Style = "Hipervínculo"
$ws.Hyperlinks.Add($ws.Range("V2"), "https://utpl.instructure.com/courses/72916")
$ws.Range("V2").Style = "Hipervínculo"

# --- View: scroll so column V is the left-most visible column ------------
$ws.Application.ActiveWindow.ScrollColumn = 22
